$wb = $excel.ActiveWorkbook

# Update the login credentials row for user 2390932 (was 2388006) with the
# refreshed password.
$wsLogin = $wb.Worksheets.Item("LoginCredentials")
$wsLogin.Range("A2").Value = 2390932
$wsLogin.Range("B2").Value = "Jul@2020"

# Refresh the NFTR ticket numbers used for escalation/filter validation.
$wsTickets = $wb.Worksheets.Item("NFTRTickets")
$wsTickets.Range("Y2").Value = "280720001032"
$wsTickets.Range("Y3").Value = "280720001033"
$wsTickets.Range("Y4").Value = "280720001034"
$wsTickets.Range("Y5").Value = "280720001035"

# Move the active tab / selection back to LoginCredentials (matches the
# author's last saved view).
$wsLogin.Activate() | Out-Null
$wsLogin.Range("D6").Select() | Out-Null
